$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column O: header "Custom Field 1", and "Test" for each data row (2-9)
$ws.Range("O1").Value = "Custom Field 1"
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 15).Value = "Test"
}

# Clear the PAN (column C) value for rows 3, 4, 7, 9 while keeping style
$ws.Range("C3").Value = $null
$ws.Range("C4").Value = $null
$ws.Range("C7").Value = $null
$ws.Range("C9").Value = $null

# Update selection to match target state
$ws.Range("O3:O9").Select()
